$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2291431.5
$ws.Range("I86").Value = 2213.5454
$ws.Range("J86").Value = 4389881.5
$ws.Range("K86").Value = 2213.5454
$ws.Range("L86").Value = 4389881.5
$ws.Range("M86").Value = -1090.5454
$ws.Range("N86").Value = -4392127.5
$ws.Range("H89").Value = 2291431.5
$ws.Range("I89").Value = 2213.5454
$ws.Range("J89").Value = 4389881.5
$ws.Range("K89").Value = 11067.727
$ws.Range("L89").Value = 21949407.5
$ws.Range("M89").Value = -5451.726999999999
$ws.Range("N89").Value = -21960639.5
$ws.Range("H100").Value = 1619.6
$ws.Range("J100").Value = 2546.75
$ws.Range("L100").Value = 2546.75
$ws.Range("N100").Value = -3628.75
$ws.Range("H103").Value = 477.375
$ws.Range("I103").Value = 345.2
$ws.Range("J103").Value = 697.6667
$ws.Range("K103").Value = 1035.6
$ws.Range("L103").Value = 2093.0001
$ws.Range("M103").Value = -449.5999999999999
$ws.Range("N103").Value = -3265.0001
$ws.Range("H125").Value = 6158.4287
$ws.Range("I125").Value = 5747.5454
$ws.Range("K125").Value = 51727.9086
$ws.Range("M125").Value = -49267.9086
$ws.Range("H129").Value = 9834.968000000001
$ws.Range("I129").Value = 29327
$ws.Range("J129").Value = 6947.2593
$ws.Range("K129").Value = 87981
$ws.Range("L129").Value = 20841.7779
$ws.Range("M129").Value = -82981
$ws.Range("N129").Value = -30841.7779
$ws.Range("H135").Value = 2300
$ws.Range("I135").Value = 858.2222
$ws.Range("J135").Value = 3597.6
$ws.Range("K135").Value = 7723.999800000001
$ws.Range("L135").Value = 32378.4
$ws.Range("M135").Value = -5188.999800000001
$ws.Range("N135").Value = -37448.39999999999
$ws.Range("H137").Value = 2291.879
$ws.Range("J137").Value = 2362.818
$ws.Range("L137").Value = 7088.454000000001
$ws.Range("N137").Value = -12188.454
$ws.Range("H138").Value = 8290.6
$ws.Range("I138").Value = 3097.5
$ws.Range("J138").Value = 8661.536
$ws.Range("K138").Value = 9292.5
$ws.Range("L138").Value = 25984.608
$ws.Range("M138").Value = -4152.5
$ws.Range("N138").Value = -36264.608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4253.222
$ws.Range("I61").Value = 1505.3636
$ws.Range("J61").Value = 8571.286
$ws.Range("K61").Value = 1505.3636
$ws.Range("L61").Value = 8571.286
$ws.Range("M61").Value = -1293.3636
$ws.Range("N61").Value = -8995.286
$ws.Range("H122").Value = 2772.372
$ws.Range("J122").Value = 3169.1292
$ws.Range("L122").Value = 9507.3876
$ws.Range("N122").Value = -14407.3876
$ws.Range("H132").Value = 237984.6
$ws.Range("I132").Value = 301902.62
$ws.Range("K132").Value = 905707.86
$ws.Range("M132").Value = -903177.86
$ws.Range("H136").Value = 4253.222
$ws.Range("I136").Value = 1505.3636
$ws.Range("J136").Value = 8571.286
$ws.Range("K136").Value = 4516.0908
$ws.Range("L136").Value = 25713.858
$ws.Range("M136").Value = -1966.0908
$ws.Range("N136").Value = -30813.858
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3656.923
$ws.Range("I99").Value = 3406
$ws.Range("K99").Value = 3406
$ws.Range("M99").Value = -1908
$ws.Range("H132").Value = 52366
$ws.Range("J132").Value = 52366
$ws.Range("L132").Value = 52366
$ws.Range("N132").Value = -62486
$ws.Range("H134").Value = 35990.582
$ws.Range("I134").Value = 3105.95
$ws.Range("K134").Value = 9317.849999999999
$ws.Range("M134").Value = -6782.849999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 811912.5600000001
$ws.Range("I31").Value = 8476026
$ws.Range("J31").Value = 5163.7896
$ws.Range("K31").Value = 8476026
$ws.Range("L31").Value = 5163.7896
$ws.Range("M31").Value = -8475731
$ws.Range("N31").Value = -5753.7896
$ws.Range("H34").Value = 811912.5600000001
$ws.Range("I34").Value = 8476026
$ws.Range("J34").Value = 5163.7896
$ws.Range("K34").Value = 8476026
$ws.Range("L34").Value = 5163.7896
$ws.Range("M34").Value = -8475824
$ws.Range("N34").Value = -5567.7896
$ws.Range("H74").Value = 99704.664
$ws.Range("J74").Value = 99704.664
$ws.Range("L74").Value = 99704.664
$ws.Range("N74").Value = -101452.664
$ws.Range("H77").Value = 99704.664
$ws.Range("J77").Value = 99704.664
$ws.Range("L77").Value = 299113.992
$ws.Range("N77").Value = -307849.992
$ws.Range("H141").Value = 90345.73
$ws.Range("J141").Value = 90345.73
$ws.Range("L141").Value = 90345.73
$ws.Range("N141").Value = -100705.73

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 6660799.5
$ws.Range("J32").Value = 825999.25
$ws.Range("L32").Value = 2477997.75
$ws.Range("N32").Value = -2478563.75
$ws.Range("H68").Value = 2438.0833
$ws.Range("J68").Value = 2543.3125
$ws.Range("L68").Value = 7629.9375
$ws.Range("N68").Value = -9251.9375
$ws.Range("H71").Value = 2438.0833
$ws.Range("J71").Value = 2543.3125
$ws.Range("L71").Value = 22889.8125
$ws.Range("N71").Value = -31001.8125
$ws.Range("H98").Value = 1767.8572
$ws.Range("J98").Value = 994
$ws.Range("L98").Value = 2982
$ws.Range("N98").Value = -5978
$ws.Range("H107").Value = 2470.641
$ws.Range("I107").Value = 918.1429000000001
$ws.Range("J107").Value = 2810.25
$ws.Range("K107").Value = 2754.4287
$ws.Range("L107").Value = 8430.75
$ws.Range("M107").Value = -834.4287000000004
$ws.Range("N107").Value = -12270.75
$ws.Range("H131").Value = 59118.113
$ws.Range("I131").Value = 84331.414
$ws.Range("K131").Value = 252994.242
$ws.Range("M131").Value = -247954.242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 57015.332
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 57015.332
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 57015.332
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -57327.332
$ws.Range("H80").Value = 731361.9399999999
$ws.Range("I80").Value = 855217
$ws.Range("J80").Value = 638470.7
$ws.Range("K80").Value = 855217
$ws.Range("L80").Value = 638470.7
$ws.Range("M80").Value = -854219
$ws.Range("N80").Value = -640466.7
$ws.Range("H83").Value = 731361.9399999999
$ws.Range("I83").Value = 855217
$ws.Range("J83").Value = 638470.7
$ws.Range("K83").Value = 4276085
$ws.Range("L83").Value = 3192353.5
$ws.Range("M83").Value = -4271093
$ws.Range("N83").Value = -3202337.5
$ws.Range("H102").Value = 6968.4614
$ws.Range("I102").Value = 8796.866
$ws.Range("K102").Value = 8796.866
$ws.Range("M102").Value = -7174.866
$ws.Range("H122").Value = 416822.34
$ws.Range("I122").Value = 533413.5600000001
$ws.Range("K122").Value = 1600240.68
$ws.Range("M122").Value = -1597790.68
$ws.Range("H132").Value = 31175.629
$ws.Range("I132").Value = 2392.2144
$ws.Range("J132").Value = 146309.28
$ws.Range("K132").Value = 7176.6432
$ws.Range("L132").Value = 438927.84
$ws.Range("M132").Value = -4646.6432
$ws.Range("N132").Value = -443987.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 505239.3
$ws.Range("I122").Value = 3477.3333
$ws.Range("K122").Value = 10431.9999
$ws.Range("M122").Value = -7981.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 16375
$ws.Range("J25").Value = 16375
$ws.Range("L25").Value = 16375
$ws.Range("N25").Value = -16961
$ws.Range("H43").Value = 26375
$ws.Range("I43").Value = 25000
$ws.Range("J43").Value = 27750
$ws.Range("K43").Value = 25000
$ws.Range("L43").Value = 27750
$ws.Range("M43").Value = -24851
$ws.Range("N43").Value = -28048
$ws.Range("H53").Value = 15250
$ws.Range("I53").Value = 14500
$ws.Range("J53").Value = 16000
$ws.Range("K53").Value = 14500
$ws.Range("L53").Value = 16000
$ws.Range("M53").Value = -13893
$ws.Range("N53").Value = -17214
$ws.Range("H55").Value = 11335.333
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 11335.333
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 11335.333
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -11889.333
$ws.Range("H107").Value = 64999.375
$ws.Range("I107").Value = 114460.445
$ws.Range("J107").Value = 1406.5714
$ws.Range("K107").Value = 343381.335
$ws.Range("L107").Value = 4219.7142
$ws.Range("M107").Value = -341461.335
$ws.Range("N107").Value = -8059.7142
$ws.Range("H132").Value = 33923.94
$ws.Range("I132").Value = 5036.7393
$ws.Range("K132").Value = 15110.2179
$ws.Range("M132").Value = -12580.2179
